$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Ch. 4) and Row 6 (Ch. 5): swap the Title / Pages values
$ws.Range("B5").Value = "Managing AD"
$ws.Range("C5").Value = 25

$ws.Range("B6").Value = "Managing Networking"
$ws.Range("C6").Value = 35

# Row 5 now has an "Achieved Delivery" date and an "Ahead or Behind" formula,
# matching the formatting already used by the other populated rows.
$ws.Range("F5").NumberFormat = $ws.Range("F4").NumberFormat
$ws.Range("F5").Value = 44690

$ws.Range("G5").Formula = "=F5-E5"
$ws.Range("G2").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the sheet's recorded selection
$ws.Range("A1:H16").Select()
